$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$Val
    )
    $rng = $ws.Range($CellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $Val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "233.70"
Set-TextValue "D3" "22.89"
Set-TextValue "D4" "5.570"
Set-TextValue "D5" "0.05667"
Set-TextValue "D6" "3.421"
Set-TextValue "D7" "6.479"
Set-TextValue "D8" "1.275"
Set-TextValue "D9" "0.8016"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D10" "0.1423"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D11" "0.07554"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D12" "0.03259"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D13" "0.02998"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D14" "0.09236"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D15" "0.001684"
$ws.Range("E15").Value = "14BitForexTokenBF"
$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D16" "3.274"
$ws.Range("E16").Value = "15MCDexMCB"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D17" "0.04740"
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D18" "0.0005990"
$ws.Range("E18").Value = "17OneONEWorstin24h"
Set-TextValue "D19" "0.006225"
Set-TextValue "D20" "0.005348"
Set-TextValue "D21" "0.001066"
Set-TextValue "D22" "0.0001504"
Set-TextValue "D23" "3.688"
Set-TextValue "D26" "0.1282"
Set-TextValue "D27" "0.0006762"
$ws.Range("E27").Value = "26UpBotsUBXTBestin24h"
Set-TextValue "D40" "0.04168"
Set-TextValue "D41" "0.006611"
Set-TextValue "D42" "0.003456"
Set-TextValue "D43" "0.1051"
Set-TextValue "D44" "0.008552"
$ws.Range("E45").Value = "44ACDXExchangeACXT"
Set-TextValue "D46" "0.00005600"
Set-TextValue "D47" "0.00000000754"
Set-TextValue "D48" "0.7892"
Set-TextValue "D49" "0.09700"
Set-TextValue "D50" "0.00002110"
Set-TextValue "D51" "0.01015"

Write-Host "Applied cryptos.xlsx symbol list update"
